# Edit script for KENTUCKY_2017.xlsx
# 1) Rename header columns to short codes
# 2) Title-case connector words (de/del/la/las/el/los/y) in state/municipality names
# 3) Remove trailing metadata rows (970-974) and shrink dimension
# 4) Minor floating point value corrections

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header renames (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2) Title-case connector words in state/municipality names ---
$ws.Range("B6").Value = "Pabellón De Arteaga"
$ws.Range("B7").Value = "Rincón De Romos"
$ws.Range("B8").Value = "San Francisco De Los Romo"
$ws.Range("B26").Value = "Amatenango De La Frontera"
$ws.Range("B29").Value = "Bejucal De Ocampo"
$ws.Range("B31").Value = "Benemérito De Las Américas"
$ws.Range("B35").Value = "Chiapa De Corzo"
$ws.Range("B38").Value = "Comitán De Domínguez"
$ws.Range("B54").Value = "Marqués De Comillas"
$ws.Range("B55").Value = "Mazapa De Madero"
$ws.Range("B61").Value = "Ocozocoautla De Espinosa"
$ws.Range("B68").Value = "San Cristóbal De Las Casas"
$ws.Range("B94").Value = "Hidalgo Del Parral"
$ws.Range("B101").Value = "San Francisco De Borja"
$ws.Range("B102").Value = "San Francisco De Conchos"
$ws.Range("B118").Value = "San Juan De Sabinas"
$ws.Range("A128").Value = "Ciudad De México"
$ws.Range("B132").Value = "Cuajimalpa De Morelos"
$ws.Range("B154").Value = "Nombre De Dios"
$ws.Range("A163").Value = "Estado De México"
$ws.Range("B163").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B165").Value = "Almoloya De Alquisiras"
$ws.Range("B166").Value = "Almoloya De Juárez"
$ws.Range("B171").Value = "Atizapán De Zaragoza"
$ws.Range("B177").Value = "Coacalco De Berriozábal"
$ws.Range("B180").Value = "Ecatepec De Morelos"
$ws.Range("B183").Value = "Ixtapan De La Sal"
$ws.Range("B184").Value = "Ixtapan Del Oro"
$ws.Range("B191").Value = "Naucalpan De Juárez"
$ws.Range("B195").Value = "San Felipe Del Progreso"
$ws.Range("B197").Value = "Soyaniquilpan De Juárez"
$ws.Range("B203").Value = "Tenango Del Valle"
$ws.Range("B208").Value = "Tlalnepantla De Baz"
$ws.Range("B211").Value = "Valle De Bravo"
$ws.Range("B212").Value = "Valle De Chalco Solidaridad"
$ws.Range("B213").Value = "Villa Del Carbón"
$ws.Range("B222").Value = "San Miguel De Allende"
$ws.Range("B223").Value = "Apaseo El Alto"
$ws.Range("B224").Value = "Apaseo El Grande"
$ws.Range("B230").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B233").Value = "Jaral Del Progreso"
$ws.Range("B239").Value = "Purísima Del Rincón"
$ws.Range("B242").Value = "San Diego De La Unión"
$ws.Range("B244").Value = "San Francisco Del Rincón"
$ws.Range("B246").Value = "San Luis De La Paz"
$ws.Range("B247").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B253").Value = "Valle De Santiago"
$ws.Range("B259").Value = "Acapulco De Juárez"
$ws.Range("B262").Value = "Ajuchitlán Del Progreso"
$ws.Range("B263").Value = "Alcozauca De Guerrero"
$ws.Range("B266").Value = "Atlamajalcingo Del Monte"
$ws.Range("B268").Value = "Atoyac De Álvarez"
$ws.Range("B269").Value = "Ayutla De Los Libres"
$ws.Range("B271").Value = "Buenavista De Cuéllar"
$ws.Range("B272").Value = "Chilapa De Álvarez"
$ws.Range("B273").Value = "Chilpancingo De Los Bravo"
$ws.Range("B274").Value = "Coyuca De Benítez"
$ws.Range("B275").Value = "Coyuca De Catalán"
$ws.Range("B278").Value = "Cutzamala De Pinzón"
$ws.Range("B282").Value = "Huitzuco De Los Figueroa"
$ws.Range("B283").Value = "Iguala De La Independencia"
$ws.Range("B284").Value = "Zihuatanejo De Azueta"
$ws.Range("B297").Value = "Taxco De Alarcón"
$ws.Range("B299").Value = "Técpan De Galeana"
$ws.Range("B301").Value = "Tepecoacuilco De Trujano"
$ws.Range("B303").Value = "Tixtla De Guerrero"
$ws.Range("B306").Value = "Tlapa De Comonfort"
$ws.Range("B316").Value = "Atotonilco El Grande"
$ws.Range("B321").Value = "Cuautepec De Hinojosa"
$ws.Range("B324").Value = "Huejutla De Reyes"
$ws.Range("B331").Value = "Mineral Del Chico"
$ws.Range("B332").Value = "Mineral Del Monte"
$ws.Range("B333").Value = "Mixquiahuala De Juárez"
$ws.Range("B334").Value = "Molango De Escamilla"
$ws.Range("B336").Value = "Pachuca De Soto"
$ws.Range("B339").Value = "Progreso De Obregón"
$ws.Range("B341").Value = "Santiago De Anaya"
$ws.Range("B342").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B345").Value = "Tenango De Doria"
$ws.Range("B346").Value = "Tepehuacán De Guerrero"
$ws.Range("B347").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B349").Value = "Tezontepec De Aldama"
$ws.Range("B354").Value = "Tula De Allende"
$ws.Range("B355").Value = "Tulancingo De Bravo"
$ws.Range("B357").Value = "Zacualtipán De Ángeles"
$ws.Range("B360").Value = "Ahualulco De Mercado"
$ws.Range("B362").Value = "Atemajac De Brizuela"
$ws.Range("B364").Value = "Atotonilco El Alto"
$ws.Range("B365").Value = "Autlán De Navarro"
$ws.Range("B373").Value = "Encarnación De Díaz"
$ws.Range("B378").Value = "Lagos De Moreno"
$ws.Range("B383").Value = "San Cristóbal De La Barranca"
$ws.Range("B384").Value = "San Juan De Los Lagos"
$ws.Range("B386").Value = "San Miguel El Alto"
$ws.Range("B389").Value = "Tepatitlán De Morelos"
$ws.Range("B391").Value = "Tlajomulco De Zúñiga"
$ws.Range("B393").Value = "Unión De San Antonio"
$ws.Range("B395").Value = "Zacoalco De Torres"
$ws.Range("B411").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B467").Value = "Coatlán Del Río"
$ws.Range("B473").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B476").Value = "Puente De Ixtla"
$ws.Range("B480").Value = "Tetela Del Volcán"
$ws.Range("B481").Value = "Tlaltizapán De Zapata"
$ws.Range("B488").Value = "Amatlán De Cañas"
$ws.Range("B490").Value = "Ixtlán Del Río"
$ws.Range("B503").Value = "Mier Y Noriega"
$ws.Range("B505").Value = "San Nicolás De Los Garza"
$ws.Range("B509").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B512").Value = "Constancia Del Rosario"
$ws.Range("B514").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B515").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B516").Value = "Ixtlán De Juárez"
$ws.Range("B517").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B521").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B522").Value = "Nejapa De Madero"
$ws.Range("B523").Value = "Oaxaca De Juárez"
$ws.Range("B565").Value = "San Pedro El Alto"
$ws.Range("B583").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B604").Value = "Santo Domingo De Morelos"
$ws.Range("B612").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B613").Value = "Tataltepec De Valdés"
$ws.Range("B614").Value = "Teotitlán De Flores Magón"
$ws.Range("B615").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B616").Value = "Tlacolula De Matamoros"
$ws.Range("B617").Value = "Villa De Tututepec"
$ws.Range("B619").Value = "Villa Sola De Vega"
$ws.Range("B620").Value = "Zapotitlán Del Río"
$ws.Range("B630").Value = "Chalchicomula De Sesma"
$ws.Range("B645").Value = "Huehuetlán El Chico"
$ws.Range("B649").Value = "Izúcar De Matamoros"
$ws.Range("B652").Value = "Los Reyes De Juárez"
$ws.Range("B653").Value = "Mazapiltepec De Juárez"
$ws.Range("B657").Value = "Palmar De Bravo"
$ws.Range("B665").Value = "Tecali De Herrera"
$ws.Range("B671").Value = "Tepexi De Rodríguez"
$ws.Range("B674").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B689").Value = "Cadereyta De Montes"
$ws.Range("B693").Value = "Jalpan De Serra"
$ws.Range("B694").Value = "Landa De Matamoros"
$ws.Range("B697").Value = "San Juan Del Río"
$ws.Range("B704").Value = "Armadillo De Los Infante"
$ws.Range("B705").Value = "Axtla De Terrazas"
$ws.Range("B707").Value = "Cerro De San Pedro"
$ws.Range("B708").Value = "Ciudad Del Maíz"
$ws.Range("B719").Value = "Santa María Del Río"
$ws.Range("B724").Value = "Villa De Arista"
$ws.Range("B725").Value = "Villa De Arriaga"
$ws.Range("B726").Value = "Villa De Guadalupe"
$ws.Range("B727").Value = "Villa De La Paz"
$ws.Range("B728").Value = "Villa De Ramos"
$ws.Range("B729").Value = "Villa De Reyes"
$ws.Range("B759").Value = "Jalpa De Méndez"
$ws.Range("B784").Value = "Soto La Marina"
$ws.Range("B802").Value = "Tepetitla De Lardizábal"
$ws.Range("B813").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B817").Value = "Amatlán De Los Reyes"
$ws.Range("B825").Value = "Boca Del Río"
$ws.Range("B830").Value = "Castillo De Teayo"
$ws.Range("B832").Value = "Cazones De Herrera"
$ws.Range("B842").Value = "Cosamaloapan De Carpio"
$ws.Range("B856").Value = "Hueyapan De Ocampo"
$ws.Range("B857").Value = "Ignacio De La Llave"
$ws.Range("B859").Value = "Ixhuatlán De Madero"
$ws.Range("B860").Value = "Ixhuatlán Del Café"
$ws.Range("B868").Value = "Juchique De Ferrer"
$ws.Range("B871").Value = "Las Vigas De Ramírez"
$ws.Range("B872").Value = "Lerdo De Tejada"
$ws.Range("B874").Value = "Martínez De La Torre"
$ws.Range("B878").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B885").Value = "Ozuluama De Mascareñas"
$ws.Range("B889").Value = "Paso De Ovejas"
$ws.Range("B890").Value = "Paso Del Macho"
$ws.Range("B893").Value = "Poza Rica De Hidalgo"
$ws.Range("B900").Value = "Sayula De Alemán"
$ws.Range("B901").Value = "Soledad De Doblado"
$ws.Range("B920").Value = "Tlacotepec De Mejía"
$ws.Range("B929").Value = "Vega De Alatorre"
$ws.Range("B943").Value = "Cañitas De Felipe Pescador"
$ws.Range("B958").Value = "Teúl De González Ortega"
$ws.Range("B961").Value = "Villa De Cos"

# --- 3) Remove trailing metadata rows 970-974 ---
$ws.Range("A970:A974").EntireRow.Delete()

# --- 4) Minor floating point corrections (last-bit rounding to match source recompute) ---
$ws.Range("D136").Value = 0.009267840593141796
$ws.Range("D326").Value = 0.009267840593141796
$ws.Range("D621").Value = 0.09128822984244672
$ws.Range("D698").Value = 0.009267840593141796
